# Horarios actualizados Línea 141 - 418
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "LP1912"
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:18:06"
$ws1.Range("A3").Value = "Total filas: 21"

# Insert two new rows before the old row 20 (pushes old rows 20-21 to 22-23)
$ws1.Rows.Item(20).Insert()
$ws1.Rows.Item(20).Insert()

$ws1.Range("A20").Value = "04:18:06"
$ws1.Range("B20").Value = "05:35"
$ws1.Range("C20").Value = "215B_EL PATO"
$ws1.Range("D20").Value = 77
$ws1.Range("E20").Value = "LP1912"

$ws1.Range("A21").Value = "04:18:06"
$ws1.Range("B21").Value = "05:36"
$ws1.Range("C21").Value = "14_ABASTO"
$ws1.Range("D21").Value = 78
$ws1.Range("E21").Value = "LP1912"

# New rows appended after the (now shifted) existing rows 22-23
$ws1.Range("A24").Value = "04:18:06"
$ws1.Range("B24").Value = "06:09"
$ws1.Range("C24").Value = "16_SANTA ANA"
$ws1.Range("D24").Value = 111
$ws1.Range("E24").Value = "LP1912"

$ws1.Range("A25").Value = "04:18:06"
$ws1.Range("B25").Value = "06:12"
$ws1.Range("C25").Value = "215A_EL PATO"
$ws1.Range("D25").Value = 114
$ws1.Range("E25").Value = "LP1912"

$ws1.Range("A26").Value = "04:18:06"
$ws1.Range("B26").Value = "06:14"
$ws1.Range("C26").Value = "225_HARAS DEL SUR"
$ws1.Range("D26").Value = 116
$ws1.Range("E26").Value = "LP1912"

# -----------------------------------------------------------------
# Sheet "LP1912-215"
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:18:06"
$ws2.Range("A3").Value = "Total filas: 8"

$ws2.Range("A12").Value = "04:18:06"
$ws2.Range("B12").Value = "05:35"
$ws2.Range("C12").Value = "215B_EL PATO"
$ws2.Range("D12").Value = 77
$ws2.Range("E12").Value = "LP1912"

$ws2.Range("A13").Value = "04:18:06"
$ws2.Range("B13").Value = "06:12"
$ws2.Range("C13").Value = "215A_EL PATO"
$ws2.Range("D13").Value = 114
$ws2.Range("E13").Value = "LP1912"

# -----------------------------------------------------------------
# Sheet "6203-6173"
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:18:06"
$ws3.Range("A3").Value = "Total filas: 4"

$ws3.Range("A8").Value = "04:18:06"
$ws3.Range("B8").Value = "05:44"
$ws3.Range("C8").Value = "215A_LA PLATA"
$ws3.Range("D8").Value = 86
$ws3.Range("E8").Value = "L6173"

$ws3.Range("A9").Value = "04:18:06"
$ws3.Range("B9").Value = "06:09"
$ws3.Range("C9").Value = "215A_LA PLATA"
$ws3.Range("D9").Value = 111
$ws3.Range("E9").Value = "L6173"
